$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("basic_amber_format")

# Fill in the new DATE / TMP rows (row 11 and row 12) that were previously blank.
$ws.Cells.Item(11, 1).Value = "DATE"
$ws.Cells.Item(11, 2).Value = "date"
$ws.Cells.Item(11, 5).Value = "object"
$ws.Cells.Item(11, 6).Value = "%Y-%m-%dT%H:%M:%S"

$ws.Cells.Item(12, 1).Value = "TMP"
$ws.Cells.Item(12, 2).Value = "temp"
$ws.Cells.Item(12, 5).Value = "float64"

# Extend the sheet with extra blank rows (13-20), matching the formatting of row 12,
# so the used range grows from A1:F12 to A1:F20.
$ws.Range("A12").Copy()
$ws.Range("A13:A20").PasteSpecial(-4122)

# Update the active selection on this sheet.
$ws.Activate()
$ws.Range("F12").Select()
